$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (row 2), pushing all
# existing data rows down by 3. Because the historical rows simply slide
# down (the window of reported weeks advances), no other manual copying
# is required - Excel's row insert naturally relocates rows 2..107 to
# rows 5..110, matching the new dimension A1:T110.
$ws.Rows("2:4").Insert()

# The insert copies formatting from the row above (the bold header row).
# Strip that back to the plain/default look used by every other data row,
# then restore the date number format on column D (style used by every
# other date cell in the column) by copying format-only from the row
# that used to be row 2 (now row 5).
$ws.Range("A2:T4").ClearFormats()
$ws.Range("D5").Copy()
$ws.Range("D2:D4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the three new rows with this week's data.
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Terminal La Palmera de La Serena"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44515
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107002
$ws.Range("J2").Value = "Chirimoya"
$ws.Range("K2").Value = "Cultivar IV Región"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 1800
$ws.Range("O2").Value = 1900
$ws.Range("P2").Value = 1850
$ws.Range("Q2").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1850
$ws.Range("T2").Value = 1

$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44515
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 1500
$ws.Range("O3").Value = 1600
$ws.Range("P3").Value = 1550
$ws.Range("Q3").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1550
$ws.Range("T3").Value = 1

$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44515
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 1200
$ws.Range("O4").Value = 1300
$ws.Range("P4").Value = 1250
$ws.Range("Q4").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1250
$ws.Range("T4").Value = 1
